# Insert a new weekly price record for "Vega Modelo de Temuco - Puerro".
# The new observation is inserted as row 291, pushing every existing row
# from 291..335 down by one (to 292..336), matching the source diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 291, shifting rows 291:335
# down to 292:336 (dimension grows from R335 to R336).
$ws.Rows.Item(291).Insert()

# Populate the newly inserted row 291 with the new record's data.
$ws.Cells.Item(291, 1).Value = 10
$ws.Cells.Item(291, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(291, 3).Value = "La Araucanía"
$ws.Cells.Item(291, 4).Value = 45180
$ws.Cells.Item(291, 5).Value = 9
$ws.Cells.Item(291, 6).Value = 100112005
$ws.Cells.Item(291, 7).Value = "Puerro"
$ws.Cells.Item(291, 8).Value = "Azul de Maquehue"
$ws.Cells.Item(291, 9).Value = "Primera"
$ws.Cells.Item(291, 10).Value = 60
$ws.Cells.Item(291, 11).Value = 9000
$ws.Cells.Item(291, 12).Value = 9000
$ws.Cells.Item(291, 13).Value = 9000
$ws.Cells.Item(291, 14).Value = "$/docena de paquetes"
$ws.Cells.Item(291, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(291, 16).Value = 750
$ws.Cells.Item(291, 17).Value = 12
$ws.Cells.Item(291, 18).Value = "Hortaliza"
